# Daily attendance processing - 2026-01-22 06:09:43
# Swap the order of names in the "Recorded By" column (G) so that entries
# reading "dnasr281@gmail.com, System" become "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    if ($cell.Value2 -eq $oldVal) {
        $cell.Value = $newVal
    }
}
